# Update the "Förändrad" (column C) date values from 45170 (2023-09-01)
# to 45174 (2023-09-05) for rows 2 through 19.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C19").Value = 45174
